$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(12).ColumnWidth = $ws.Columns.Item(12).ColumnWidth
